$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells we touch keep their original text formatting
# (avoids Excel auto-converting numeric-looking strings like "19.26" into numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.771.70"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.648.09"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.251"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "19.26"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "1.878.26"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "1.632.06"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "0.531"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "66.09"
$ws.Range("E16").Value = "  +4.86%  "
$ws.Range("D17").Value = "26.819.10"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "0.0₃0746"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "218.29"
$ws.Range("E19").Value = "  +4.08%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("E22").Value = "  +2.94%  "
$ws.Range("D23").Value = "9.46"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  +9.22%  "
$ws.Range("D25").Value = "147.85"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "6.94"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "15.79"
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  +4.22%  "
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("D34").Value = "1.272.49"
$ws.Range("E34").Value = "  +9.15%  "
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("E37").Value = "  +4.46%  "
$ws.Range("D38").Value = "0.808"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").Value = "0.514"
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("D42").Value = "0.807"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D43").Value = "5.37"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "1.787.46"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("D45").Value = "93.71"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("E46").Value = "  +4.26%  "
$ws.Range("D47").Value = "55.96"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "7.68"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("D51").Value = "0.0968"
$ws.Range("E51").Value = "  +3.38%  "
